# Updated cryptos list — refresh prices / 1h volume deltas, and fix the
# rank ordering for two coin pairs (WrappedEther/Polkadot and
# Stellar/ImmutableX swapped places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, regardless of whether it "looks"
# numeric (e.g. "330.93", "1.000"), matching the existing inline-string
# cells, and without leaving a lingering NumberFormat/quote-prefix style
# on the cell (reset to the default "Normal" style after writing).
function Set-TextCell {
    param($sheet, $addr, $val)
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# --- Row 13 / Row 14: WrappedEther and Polkadot swap rank order ---------
Set-TextCell $ws "B13" "WrappedEther"
Set-TextCell $ws "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws "D13" "1.835.10"
Set-TextCell $ws "E13" "  -1.77%  "

Set-TextCell $ws "B14" "Polkadot"
Set-TextCell $ws "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D14" "5.880"
Set-TextCell $ws "E14" "  -2.33%  "

# --- Row 32 / Row 33: Stellar and ImmutableX swap rank order -------------
Set-TextCell $ws "B32" "Stellar"
Set-TextCell $ws "C32" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D32" "0.09316"
Set-TextCell $ws "E32" "  -2.36%  "

Set-TextCell $ws "B33" "ImmutableX"
Set-TextCell $ws "C33" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws "D33" "0.9388"
Set-TextCell $ws "E33" "  -5.13%  "

# --- Remaining rows: price (D) / volume-1h (E) refresh --------------------
Set-TextCell $ws "D2" "27.428.60"
Set-TextCell $ws "E2" "  -1.57%  "

Set-TextCell $ws "D3" "1.828.04"
Set-TextCell $ws "E3" "  -2.80%  "

Set-TextCell $ws "D5" "330.93"
Set-TextCell $ws "E5" "  -1.57%  "

Set-TextCell $ws "D6" "1.002"
Set-TextCell $ws "E6" "  -0.64%  "

Set-TextCell $ws "D7" "0.4579"
Set-TextCell $ws "E7" "  -2.68%  "

Set-TextCell $ws "D8" "0.3820"
Set-TextCell $ws "E8" "  -3.33%  "

Set-TextCell $ws "D9" "46.28"
Set-TextCell $ws "E9" "  +1.48%  "

Set-TextCell $ws "D10" "0.07893"
Set-TextCell $ws "E10" "  -1.78%  "

Set-TextCell $ws "D11" "0.9662"
Set-TextCell $ws "E11" "  -4.79%  "

Set-TextCell $ws "E12" "  -4.57%  "

Set-TextCell $ws "D15" "7.092"
Set-TextCell $ws "E15" "  -2.90%  "

Set-TextCell $ws "D16" "1.002"
Set-TextCell $ws "E16" "  -0.72%  "

Set-TextCell $ws "D17" "89.80"
Set-TextCell $ws "E17" "  +0.66%  "

Set-TextCell $ws "D18" "0.06611"
Set-TextCell $ws "E18" "  -1.72%  "

Set-TextCell $ws "E19" "  -2.09%  "

Set-TextCell $ws "D20" "17.21"
Set-TextCell $ws "E20" "  -0.89%  "

Set-TextCell $ws "D21" "1.000"
Set-TextCell $ws "E21" "  -1.08%  "

Set-TextCell $ws "D22" "27.423.74"
Set-TextCell $ws "E22" "  -1.51%  "

Set-TextCell $ws "D23" "5.315"
Set-TextCell $ws "E23" "  -3.47%  "

Set-TextCell $ws "D24" "10.85"
Set-TextCell $ws "E24" "  -1.60%  "

Set-TextCell $ws "D25" "2.280"
Set-TextCell $ws "E25" "  -1.58%  "

Set-TextCell $ws "D26" "2.045.41"
Set-TextCell $ws "E26" "  -2.36%  "

Set-TextCell $ws "D27" "156.38"
Set-TextCell $ws "E27" "  -1.97%  "

Set-TextCell $ws "D28" "19.41"
Set-TextCell $ws "E28" "  -2.45%  "

Set-TextCell $ws "D29" "2.070"
Set-TextCell $ws "E29" "  -4.52%  "

Set-TextCell $ws "D30" "5.281"
Set-TextCell $ws "E30" "  -4.04%  "

Set-TextCell $ws "D31" "118.19"
Set-TextCell $ws "E31" "  -3.28%  "

Set-TextCell $ws "D34" "3.576"
Set-TextCell $ws "E34" "  -1.59%  "

Set-TextCell $ws "D35" "5.231"
Set-TextCell $ws "E35" "  -2.36%  "

Set-TextCell $ws "D36" "1.326"
Set-TextCell $ws "E36" "  -2.18%  "

Set-TextCell $ws "D37" "0.05937"
Set-TextCell $ws "E37" "  -2.42%  "

Set-TextCell $ws "D38" "0.02181"
Set-TextCell $ws "E38" "  -3.00%  "

Set-TextCell $ws "D39" "8.140"
Set-TextCell $ws "E39" "  -2.45%  "

Set-TextCell $ws "D40" "1.001"
Set-TextCell $ws "E40" "  -0.64%  "

Set-TextCell $ws "D41" "1.145"
Set-TextCell $ws "E41" "  -4.75%  "

Set-TextCell $ws "D42" "0.5779"
Set-TextCell $ws "E42" "  -3.93%  "

Set-TextCell $ws "D43" "0.1825"
Set-TextCell $ws "E43" "  -3.83%  "

Set-TextCell $ws "D44" "9.996"
Set-TextCell $ws "E44" "  -4.52%  "

Set-TextCell $ws "D45" "1.267"
Set-TextCell $ws "E45" "  +1.51%  "

Set-TextCell $ws "D46" "11.96"
Set-TextCell $ws "E46" "  -2.26%  "

Set-TextCell $ws "D47" "0.5429"
Set-TextCell $ws "E47" "  -4.65%  "

Set-TextCell $ws "D48" "1.874"
Set-TextCell $ws "E48" "  -3.69%  "

Set-TextCell $ws "D50" "109.64"
Set-TextCell $ws "E50" "  -2.71%  "

Set-TextCell $ws "D51" "1.002"
Set-TextCell $ws "E51" "  -33.95%  "
